$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 380. This shifts the existing rows
# 380..502 down to 383..505 (preserving all of their values/formats),
# matching the diff where every record from the old row 380 onward moves
# down by exactly 3 rows.
$ws.Rows("380:382").Insert()

# Populate the 3 newly-inserted rows (380..382) with the new weekly
# record (same constant columns as every other row in this sheet).
$newRows = @(
    @{ Row = 380; D = 44588; K = "Sin especificar"; L = "Pintón";          M = 80;  N = 14000; S = 700 },
    @{ Row = 381; D = 44588; K = "Sin especificar"; L = "Primera Maduro";  M = 120; N = 16000; S = 800 },
    @{ Row = 382; D = 44588; K = "Sin especificar"; L = "Primera Pintón"; M = 120; N = 17000; S = 850 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = 8
    $ws.Cells.Item($row, 2).Value2  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value2  = "Coquimbo"
    $ws.Cells.Item($row, 4).Value2  = $r.D
    $ws.Cells.Item($row, 5).Value2  = 4
    $ws.Cells.Item($row, 6).Value2  = "Fruta"
    $ws.Cells.Item($row, 7).Value2  = 100108
    $ws.Cells.Item($row, 8).Value2  = "Tropicales y subtropicales"
    $ws.Cells.Item($row, 9).Value2  = 100108006
    $ws.Cells.Item($row, 10).Value2 = "Plátano"
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.N
    $ws.Cells.Item($row, 16).Value2 = $r.N
    $ws.Cells.Item($row, 17).Value2 = "$/caja 20 kilos"
    $ws.Cells.Item($row, 18).Value2 = "Ecuador"
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = 20
}

# Append 3 new rows at the end (503..505), duplicating the record that
# used to be the final block (now at 500..502 after the shift above).
$lastRows = @(
    @{ Row = 503; D = 44392; K = "Sin especificar"; L = "Pintón";          M = 80;  N = 13000; S = 650 },
    @{ Row = 504; D = 44392; K = "Sin especificar"; L = "Primera Maduro";  M = 120; N = 14500; S = 725 },
    @{ Row = 505; D = 44392; K = "Sin especificar"; L = "Primera Pintón"; M = 120; N = 15000; S = 750 }
)

foreach ($r in $lastRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = 8
    $ws.Cells.Item($row, 2).Value2  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value2  = "Coquimbo"
    $ws.Cells.Item($row, 4).Value2  = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value2  = 4
    $ws.Cells.Item($row, 6).Value2  = "Fruta"
    $ws.Cells.Item($row, 7).Value2  = 100108
    $ws.Cells.Item($row, 8).Value2  = "Tropicales y subtropicales"
    $ws.Cells.Item($row, 9).Value2  = 100108006
    $ws.Cells.Item($row, 10).Value2 = "Plátano"
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.N
    $ws.Cells.Item($row, 16).Value2 = $r.N
    $ws.Cells.Item($row, 17).Value2 = "$/caja 20 kilos"
    $ws.Cells.Item($row, 18).Value2 = "Ecuador"
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = 20
}
